$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @("2025-08-25", "[]"),
    @("2025-08-26", "[]"),
    @("2025-08-27", "[]"),
    @("2025-08-28", "[]")
)

$startRow = 75
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $newRows[$i][0]
    $cellA.Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
}
